# Pooh Points: final 20260211 -> PD12
# Mark remaining in-progress games as "Final", adjust one player's minutes,
# and tighten the "status" column width now that it only needs to fit "Final".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Players")

# Rows in the "status" column (G) that are not yet "Final" and need updating.
$rowsToFinalize = @(6,7,8,13,15,17,21,24,32,33,34,40,43,45,54,56,57,78,79,80,82,83,85,86,89,90,92,93,95,100,101,102,106,107,113,114,117)

foreach ($r in $rowsToFinalize) {
    $ws.Cells.Item($r, 7).Value = "Final"
}

# Jayden Stone's minutes (row 32, column P = min) corrected from 30 to 29.
$ws.Cells.Item(32, 16).Value = 29

# Column G ("status") no longer needs to fit long clock strings; narrow it
# from width 17 to width 8 now that every value is the short word "Final".
$ws.Columns(7).ColumnWidth = 7.17
